# Apply "new quests for the new plane" changes to the Locations sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Existing row 14: set the "type" column (L) to 3.
$ws.Range("L14").Value = 3

# New row 44: Twisted Dimensional Gate (Hell plane)
$ws.Range("A44").Value = "Twisted Dimensional Gate"
$ws.Range("B44").Value = "Hell"
$ws.Range("D44").Value = "Twisted Tree Branch"
$ws.Range("E44").Value = "A place that leads to the delisions of those who live in a fantasy world. Unable to esxape their past, unwilling to let go. The world that once was now becomes twisted and broken in the eyes of The Wondering Prince"
$ws.Range("G44").Value = 1
$ws.Range("J44").Value = 464
$ws.Range("K44").Value = 64
$ws.Range("L44").Value = 6
$ws.Range("M44").Value = "No"

$wb.Save()
